# Binary glTF figures: widen header 16-byte -> 20-byte, insert a new
# "length" field between "version" and "jsonOffset", and shift every
# shape to the right of it (plus a small uniform +6350 EMU vertical
# nudge on the header row) to keep the diagram aligned.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU -> points, with a +0.5 EMU bias that compensates for the
# single-precision float truncation the COM layer applies when it
# converts Shape.Left/Top/Width/Height (points) back to EMU on save.
function EMU([double]$emu) {
    return ($emu + 0.5) / 12700.0
}

function SetPos($shp, $left, $top, $width, $height) {
    if ($null -ne $left)   { $shp.Left   = EMU $left }
    if ($null -ne $top)    { $shp.Top    = EMU $top }
    if ($null -ne $width)  { $shp.Width  = EMU $width }
    if ($null -ne $height) { $shp.Height = EMU $height }
}

function ById($id) {
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.Id -eq $id) { return $shp }
    }
    return $null
}

# -- magic --------------------------------------------------------------
SetPos (ById 4) $null 662695 $null $null

# -- version --------------------------------------------------------------
SetPos (ById 5) $null 662695 $null $null

# -- jsonOffset (shifts right to make room for the new "length" field) --
SetPos (ById 6) 3813235 662695 $null $null

# -- jsonLength --
SetPos (ById 7) 4921381 662695 $null $null

# -- JSON --
SetPos (ById 8) 6952977 662695 $null $null

# -- embedded data (dashed, left of JSON) --
SetPos (ById 9) 6029527 662695 $null $null

# -- embedded data (dashed, right of JSON) --
SetPos (ById 10) 7784079 662695 $null $null

# -- External data oval --
SetPos (ById 11) 6701051 1888153 $null $null

# -- header/body top bracket arrow (widens) --
SetPos (ById 13) $null $null 5919720 $null

# -- "16-byte header" -> "20-byte header" label, split across two runs --
$hdrShape = ById 14
SetPos $hdrShape 2330898 $null $null $null
$hdrRange = $hdrShape.TextFrame.TextRange
$hdrRange.Text = "20-"
$hdrRange.InsertAfter("byte header") | Out-Null

# -- body bracket arrow --
SetPos (ById 15) 6029527 $null $null $null

# -- "body" label --
SetPos (ById 17) 7090247 $null $null $null

# -- curved connector (JSON -> embedded data, right) --
SetPos (ById 23) 7807166 692330 $null $null

# -- curved connector (JSON -> embedded data, left) --
SetPos (ById 24) 6929890 692330 $null $null

# -- arrow connector down to "External data" oval --
SetPos (ById 35) 7368528 1137318 3810 750835

# -- "bufferViews" label --
SetPos (ById 36) 7335020 1382048 $null $null

# -- new "length (uint32)" field, inserted where jsonOffset used to be --
# Consume the shape ids that would otherwise be handed out first (2, 3,
# 12, 16) with throw-away textboxes so the real new shape lands on id 18,
# matching the id the rest of the deck already expects.
$placeholderIds = 2, 3, 12, 16
$throwaways = @()
foreach ($unused in $placeholderIds) {
    $throwaways += $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
}
foreach ($t in $throwaways) { $t.Delete() }

$versionShape = ById 5
$lengthRange = $versionShape.Duplicate()
$lengthShape = $lengthRange.Item(1)
$lengthShape.Name = "TextBox 17"
SetPos $lengthShape 2880187 662695 923450 461665

# First paragraph of the duplicated shape reads "version" (7 chars);
# replace it with "length" while keeping its run formatting untouched.
$lengthShape.TextFrame.TextRange.Characters(1, 7).Text = "length"
